# Rename transcript speaker labels in the DataSheet:
#   "Davis"   -> "T"
#   "Student" -> "S"
# Only cells in column D (Speaker) are affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)   # Column D
    $val = $cell.Value2
    if ($val -eq "Davis") {
        $cell.Value2 = "T"
    } elseif ($val -eq "Student") {
        $cell.Value2 = "S"
    }
}
